# Update comparison plots, thicker traces, list of symbols
# Applies the numeric/style/selection changes to Sheet1 that back the four
# scatter charts on the sheet. Excel recalculates the chart caches from the
# cell values on save, so editing the backing cells is sufficient to move
# the charts' plotted points.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Stage 1 / default) -------------------------------------------
$ws.Cells.Item(4, 7).Value = 0.67    # G4
$ws.Cells.Item(4, 8).Value = 0.77    # H4
$ws.Cells.Item(4, 9).Value = 0.35    # I4
$ws.Cells.Item(4, 10).Value = 0.69   # J4

# --- Row 6 (Stage 1 / experiment) -----------------------------------------
# Values change and the stale yellow "to fill in" highlight is cleared
# (equivalent to the Excel "No Fill" command) now that the numbers are in.
$ws.Cells.Item(6, 9).Value = 0.62    # I6
$ws.Cells.Item(6, 9).Interior.ColorIndex = -4142   # xlColorIndexNone
$ws.Cells.Item(6, 10).Value = 0.78   # J6
$ws.Cells.Item(6, 10).Interior.ColorIndex = -4142  # xlColorIndexNone

# --- Row 7 (Stage 3 / default) --------------------------------------------
# Previously blank placeholder cells (just highlighted yellow); now filled
# in with data, with the highlight cleared.
$rowData7 = @{ 3 = 0.58; 4 = 0.58; 5 = 0.33; 6 = 0.58; 7 = 0.7; 8 = 0.75; 9 = 0.38; 10 = 0.77 }
foreach ($col in $rowData7.Keys) {
    $cell = $ws.Cells.Item(7, $col)
    $cell.Value = $rowData7[$col]
    $cell.Interior.ColorIndex = -4142   # xlColorIndexNone
}

# --- Row 11 (Stage 4 / simulation) ----------------------------------------
# Values are unchanged; the explicit "no fill" formatting on I11/J11 is
# reset back to the plain Normal cell style.
$ws.Cells.Item(11, 9).Style = "Normal"
$ws.Cells.Item(11, 10).Style = "Normal"

# --- Selection --------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("I9:J9").Select() | Out-Null

# --- Make sure everything (including the charts) is recalculated ---------
$co = $ws.ChartObjects()
for ($i = 1; $i -le $co.Count; $i++) {
    try { $co.Item($i).Chart.Refresh() | Out-Null } catch {}
}
$wb.RefreshAll() | Out-Null
$excel.CalculateFullRebuild() | Out-Null
